$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 171: Add / הוסף  (B171 reuses the existing "Add" shared string, same
# style as the preceding key cells; C171 is a brand new string, default style)
$ws.Range("B171").Value = "Add"
$ws.Range("C171").Value = "הוסף"

# Row 172: Payment status / מצב תשלום  (same key-cell style as row 171)
$ws.Range("B172").Value = "Payment status"
$ws.Range("C172").Value = "מצב תשלום"

# Row 173: unfulfilled / לא שולם  (key cell uses the alternate "code" style)
$ws.Range("B173").Value = "unfulfilled"
$ws.Range("C173").Value = "לא שולם"

# Row 174: fulfilled / שולם  (key cell uses the plain default style)
$ws.Range("B174").Value = "fulfilled"
$ws.Range("C174").Value = "שולם"

# Copy just the cell formatting (not values) from existing cells that carry
# the same styles the diff shows for the new rows.
$ws.Range("B170").Copy()
$ws.Range("B171:B172").PasteSpecial(-4122)

$ws.Range("B132").Copy()
$ws.Range("B173").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection to B174, matching the final saved view state.
$ws.Range("B174").Select()
